# Update the Version Control table: reassign reviewer/responsible names
# for the first three revision rows (rows 6, 7, 8 of the table, 1-based).
#
# Row 6 (ver 2.8.1):  วรรัตน์ (QM)  -> ณัฐดนัย (DM)
#                     กิตติพศ (SP) -> วิรัตน์ (TL)
# Row 7 (ver 1.5.1):  วรรัตน์ (QM)  -> วิรัตน์ (TL)
#                     กิตติพศ (SP) -> unchanged
# Row 8 (ver 1.4.2):  วรรัตน์ (QM)  -> ณัฐนันท์ (QA)
#                     กิตติพศ (SP) -> unchanged

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-InCell($cell, $oldText, $newText, $wildcards) {
    $rng = $cell.Range
    # Replace:=1 (wdReplaceOne) -- each search string occurs exactly once
    # within the targeted cell, and ReplaceAll(=2) in this COM host is not
    # confined to the supplied Range, so it must be avoided here.
    $rng.Find.Execute($oldText, $false, $false, $wildcards, $false, $false, $true, 1, $false, $newText, 1) | Out-Null
}

# --- Row 6: "ผู้รับผิดชอบ" cell (column 4) ---
$cell = $t.Rows.Item(6).Cells.Item(4)
Replace-InCell $cell "วรรัตน์ " "ณัฐดนัย" $false
Replace-InCell $cell "(QM)" " (DM)" $false

# --- Row 6: "ผู้ตรวจ" cell (column 5) ---
$cell = $t.Rows.Item(6).Cells.Item(5)
Replace-InCell $cell "กิตติพศ " "วิรัตน์" $false
Replace-InCell $cell "(SP)" " (TL)" $false

# --- Row 7: "ผู้รับผิดชอบ" cell (column 4) ---
$cell = $t.Rows.Item(7).Cells.Item(4)
Replace-InCell $cell "วรรัตน์ " "วิรัตน์" $false
Replace-InCell $cell "(QM)" " (TL)" $false

# --- Row 8: "ผู้รับผิดชอบ" cell (column 4) ---
$cell = $t.Rows.Item(8).Cells.Item(4)
Replace-InCell $cell "วรรัตน์ " "ณัฐนันท์ " $false
Replace-InCell $cell "(QM)" "(QA)" $false
